$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text that can look numeric (e.g.
# "1.00", "34.446.33"). Before writing each such cell, force its number
# format to Text ("@") so Excel keeps the literal string instead of
# silently reinterpreting/renormalising it as a number. Only the specific
# cells being rewritten are touched; everything else is left alone.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.446.33"
$ws.Range("E2").Value = "  -0.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.34"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.05"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("E6").Value = "  +3.16%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.11"
$ws.Range("E8").Value = "  +6.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.286"
$ws.Range("E9").Value = "  -4.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0664"
$ws.Range("E10").Value = "  -4.27%  "

$ws.Range("E11").Value = "  +1.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.061.67"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.83"
$ws.Range("E13").Value = "  -6.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.795.72"
$ws.Range("E14").Value = "  -1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.426.32"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "0.627"
$ws.Range("E16").Value = "  -4.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  -3.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.83"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.14"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0763"
$ws.Range("E20").Value = "  -3.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.06"
$ws.Range("E21").Value = "  -4.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("E23").Value = "  -3.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -3.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.21"
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.51"
$ws.Range("E26").Value = "  +3.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.61"
$ws.Range("E27").Value = "  -4.57%  "

$ws.Range("E28").Value = "  +2.07%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.73"
$ws.Range("E31").Value = "  -3.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0511"
$ws.Range("E32").Value = "  -3.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.81"
$ws.Range("E33").Value = "  -5.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -0.83%  "

$ws.Range("E35").Value = "  -0.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.298.95"
$ws.Range("E36").Value = "  -6.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.632"
$ws.Range("E37").Value = "  -6.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0184"
$ws.Range("E38").Value = "  -3.27%  "

$ws.Range("E39").Value = "  +0.29%  "

$ws.Range("E40").Value = "  -6.92%  "

$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "81.31"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("E43").Value = "  -1.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.939"
$ws.Range("E44").Value = "  -3.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("E46").Value = "  +3.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.962.36"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.71"
$ws.Range("E48").Value = "  -5.39%  "

$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.18"
$ws.Range("E50").Value = "  -3.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0611"
$ws.Range("E51").Value = "  -0.96%  "
